$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as literal text,
# even when the text looks like a plain number (avoids Excel auto-numeric conversion).
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    if ($val -match "^[+-]?[0-9]*\.?[0-9]+$") {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $val
}

Set-TextValue "D2" '41.433.88'
Set-TextValue "E2" '  -0.56%  '
Set-TextValue "D3" '2.443.42'
Set-TextValue "E3" '  -1.23%  '
Set-TextValue "E4" '  +0.27%  '
Set-TextValue "D5" '317.08'
Set-TextValue "E5" '  -0.38%  '
Set-TextValue "D6" '90.27'
Set-TextValue "E6" '  -2.39%  '
Set-TextValue "D7" '0.544'
Set-TextValue "E7" '  -1.87%  '
Set-TextValue "E8" '  +0.13%  '
Set-TextValue "D9" '0.499'
Set-TextValue "E9" '  -3.18%  '
Set-TextValue "D10" '32.28'
Set-TextValue "E10" '  -1.60%  '
Set-TextValue "D11" '0.0832'
Set-TextValue "E11" '  -6.80%  '
Set-TextValue "D12" '0.109'
Set-TextValue "E12" '  -1.32%  '
Set-TextValue "D13" '2.820.97'
Set-TextValue "E13" '  -1.11%  '
Set-TextValue "D14" '6.72'
Set-TextValue "E14" '  -2.47%  '
Set-TextValue "D15" '15.38'
Set-TextValue "E15" '  -1.90%  '
Set-TextValue "D16" '2.430.48'
Set-TextValue "E16" '  -0.88%  '
Set-TextValue "E17" '  -1.70%  '
Set-TextValue "D18" '41.338.22'
Set-TextValue "D19" '6.28'
Set-TextValue "E19" '  -2.79%  '
Set-TextValue "E20" '  -4.10%  '
Set-TextValue "D21" '72.08'
Set-TextValue "E21" '  +0.81%  '
Set-TextValue "D22" '11.16'
Set-TextValue "E22" '  -2.54%  '
Set-TextValue "D23" '235.17'
Set-TextValue "E23" '  -2.67%  '
Set-TextValue "D24" '2.70'
Set-TextValue "E24" '  -1.60%  '
Set-TextValue "E25" '  +0.10%  '
Set-TextValue "D26" '1.89'
Set-TextValue "E26" '  -1.43%  '
Set-TextValue "D27" '24.12'
Set-TextValue "E27" '  -3.13%  '
Set-TextValue "D28" '2.22'
Set-TextValue "E28" '  -2.83%  '
Set-TextValue "D29" '9.58'
Set-TextValue "E29" '  -2.63%  '
Set-TextValue "D30" '34.91'
Set-TextValue "E30" '  -2.61%  '
Set-TextValue "D31" '157.12'
Set-TextValue "E31" '  +0.38%  '
Set-TextValue "E32" '  -3.90%  '
Set-TextValue "E33" '  +0.09%  '
Set-TextValue "E34" '  -1.32%  '
Set-TextValue "D35" '0.0747'
Set-TextValue "E35" '  -2.51%  '
Set-TextValue "E36" '  +0.18%  '
Set-TextValue "D37" '16.64'
Set-TextValue "E37" '  -4.49%  '
Set-TextValue "E38" '  -0.74%  '
Set-TextValue "E39" '  -2.43%  '
Set-TextValue "D40" '0.100'
Set-TextValue "E40" '  -2.14%  '
Set-TextValue "D41" '3.89'
Set-TextValue "E41" '  -2.08%  '
Set-TextValue "E42" '  -7.19%  '
Set-TextValue "D43" '1.987.28'
Set-TextValue "E43" '  +0.48%  '
Set-TextValue "E44" '  -2.98%  '
Set-TextValue "D45" '18.15'
Set-TextValue "E45" '  -5.55%  '
Set-TextValue "E46" '  -3.54%  '
Set-TextValue "D47" '9.51'
Set-TextValue "E47" '  +4.13%  '
Set-TextValue "D48" '2.679.28'
Set-TextValue "E48" '  -1.00%  '
Set-TextValue "D49" '95.54'
Set-TextValue "E49" '  -1.97%  '
Set-TextValue "D50" '73.50'
Set-TextValue "E50" '  -0.77%  '

# Row 51 full update (ordi -> MultiversX)
$ws.Range("B51").Value = "MultiversX"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
Set-TextValue "D51" '52.20'
Set-TextValue "E51" '  -0.68%  '

Write-Host "Applied cryptos update"
